$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

# E2: set new text value
$ws.Range("E2").Value = "20 TL - 20 TL"

# F3, F4, F5: clear existing text (becomes empty inlineStr)
$ws.Range("F3").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("F5").Value = ""

# F8, F9, F10: clear existing text
$ws.Range("F8").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("F10").Value = ""

# E13: update text value
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 8.700 TL"

# F13: clear existing text
$ws.Range("F13").Value = ""

# F14: clear existing text
$ws.Range("F14").Value = ""
